$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking strings
# (e.g. "499.42") are stored as text, matching the source data (t="inlineStr").
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '58.318.05'
$ws.Range("E2").Value = '  -3.30%  '
$ws.Range("D3").Value = '2.689.85'
$ws.Range("E3").Value = '  -7.47%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '499.42'
$ws.Range("E5").Value = '  -5.40%  '
$ws.Range("D6").Value = '138.68'
$ws.Range("E6").Value = '  -2.48%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").Value = '0.525'
$ws.Range("E8").Value = '  -4.77%  '
$ws.Range("D9").Value = '2.691.74'
$ws.Range("E9").Value = '  -7.40%  '
$ws.Range("E10").Value = '  +1.72%  '
$ws.Range("E11").Value = '  -5.22%  '
$ws.Range("D12").Value = '0.341'
$ws.Range("E12").Value = '  -3.34%  '
$ws.Range("E13").Value = '  +1.16%  '
$ws.Range("D14").Value = '3.156.58'
$ws.Range("E14").Value = '  -7.54%  '
$ws.Range("D15").Value = '58.279.65'
$ws.Range("E15").Value = '  -3.60%  '
$ws.Range("D16").Value = '21.22'
$ws.Range("E16").Value = '  -6.29%  '
$ws.Range("D17").Value = '2.691.16'
$ws.Range("E17").Value = '  -7.29%  '
$ws.Range("E18").Value = '  -5.46%  '
$ws.Range("D19").Value = '4.64'
$ws.Range("E19").Value = '  -5.73%  '
$ws.Range("D20").Value = '10.77'
$ws.Range("E20").Value = '  -6.78%  '
$ws.Range("D21").Value = '336.19'
$ws.Range("E21").Value = '  -6.74%  '
$ws.Range("D22").Value = '6.11'
$ws.Range("E22").Value = '  -7.47%  '
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("E24").Value = '  -0.53%  '
$ws.Range("D25").Value = '61.86'
$ws.Range("E25").Value = '  -2.16%  '
$ws.Range("B26").Value = 'Polygon'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D26").Value = '0.418'
$ws.Range("E26").Value = '  -7.03%  '
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").Value = '0.169'
$ws.Range("E27").Value = '  -3.86%  '
$ws.Range("D28").Value = '0.996'
$ws.Range("E28").Value = '  -0.59%  '
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").Value = '7.33'
$ws.Range("E29").Value = '  -4.68%  '
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0812'
$ws.Range("E30").Value = '  -5.08%  '
$ws.Range("E31").Value = '  -0.15%  '
$ws.Range("E32").Value = '  -4.93%  '
$ws.Range("D33").Value = '18.80'
$ws.Range("E33").Value = '  -4.46%  '
$ws.Range("D34").Value = '147.43'
$ws.Range("E34").Value = '  -4.05%  '
$ws.Range("D35").Value = '4.10'
$ws.Range("E35").Value = '  -4.64%  '
$ws.Range("D36").Value = '5.26'
$ws.Range("E36").Value = '  -5.10%  '
$ws.Range("D37").Value = '0.915'
$ws.Range("E37").Value = '  -8.24%  '
$ws.Range("E38").Value = '  -6.85%  '
$ws.Range("D39").Value = '35.82'
$ws.Range("E39").Value = '  -5.30%  '
$ws.Range("D40").Value = '1.37'
$ws.Range("E40").Value = '  -6.13%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '2.146.91'
$ws.Range("E41").Value = '  -8.18%  '
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").Value = '0.996'
$ws.Range("E42").Value = '  -0.22%  '
$ws.Range("D43").Value = '3.46'
$ws.Range("E43").Value = '  -5.61%  '
$ws.Range("D44").Value = '0.0549'
$ws.Range("E44").Value = '  -3.48%  '
$ws.Range("D45").Value = '0.589'
$ws.Range("E45").Value = '  -8.30%  '
$ws.Range("D46").Value = '10.33'
$ws.Range("E46").Value = '  -0.22%  '
$ws.Range("D47").Value = '18.49'
$ws.Range("E47").Value = '  -11.39%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").Value = '0.0223'
$ws.Range("E48").Value = '  -4.17%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '4.57'
$ws.Range("E49").Value = '  -5.85%  '
$ws.Range("D50").Value = '0.0877'
$ws.Range("E50").Value = '  -5.11%  '
$ws.Range("D51").Value = '17.49'
$ws.Range("E51").Value = '  -4.43%  '

# Restore the default (unstyled) cell style so the style index matches the original.
$dataRange.Style = "Normal"
